$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26 (shifts FATURADO..TRANSFERÊNCIA down by one)
$ws.Rows.Item(26).Insert()

# Set the new cell's value
$ws.Range("A26").Value = "JUROS"

# Update the defined name range to cover the new row count (A1:A37)
$wb.Names.Item("Forma_de_Pagamento").RefersTo = "='Forma_de_Pagamento'!`$A`$1:`$A`$37"
